$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are treated as text so values like "36.504.77" or "242.88"
# are not auto-converted to numbers/dates by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.504.77'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.951.00'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.88'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.615'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.96'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +4.48%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.374'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0786'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -7.00%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.79%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.15'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +5.29%  '
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.242.48'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.62%  '
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.827'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.49'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.22'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.94%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.952.66'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.419.13'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.33'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0847'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '228.62'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.05'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.13%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.32%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.68%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +6.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.14'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '159.69'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.20'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.30'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +18.48%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.72'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0609'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.42'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +5.42%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.19%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.44'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +8.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.26'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.80%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.43'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -12.85%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0949'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.46%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.16'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0208'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.73'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.357.82'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.41'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.02'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.09'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.43'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +4.91%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.132.33'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.56%  '
